$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on cells whose new values would otherwise be
# auto-parsed by Excel as numbers (single-dot decimal-looking strings),
# so they are preserved as literal text, matching the source data feed.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values (coin names, links, prices, volume deltas).
$ws.Range("D2").Value = "27.207.49"
$ws.Range("E2").Value = "  -1.33%  "
$ws.Range("D3").Value = "1.786.72"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("E4").Value = "  +0.44%  "
$ws.Range("D5").Value = "335.94"
$ws.Range("E5").Value = "  -2.63%  "
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("D7").Value = "0.3880"
$ws.Range("E7").Value = "  +1.15%  "
$ws.Range("D8").Value = "0.3424"
$ws.Range("E8").Value = "  -2.75%  "
$ws.Range("D9").Value = "48.38"
$ws.Range("E9").Value = "  -3.47%  "
$ws.Range("D10").Value = "1.197"
$ws.Range("E10").Value = "  -3.13%  "
$ws.Range("D11").Value = "0.07503"
$ws.Range("E11").Value = "  -3.11%  "
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").Value = "21.84"
$ws.Range("E13").Value = "  -3.12%  "
$ws.Range("D14").Value = "6.465"
$ws.Range("E14").Value = "  -2.42%  "
$ws.Range("D15").Value = "1.790.95"
$ws.Range("E15").Value = "  -1.33%  "
$ws.Range("D16").Value = "7.090"
$ws.Range("E16").Value = "  -1.67%  "
$ws.Range("E17").Value = "  -2.80%  "
$ws.Range("D18").Value = "0.06680"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("D19").Value = "83.88"
$ws.Range("E19").Value = "  -3.29%  "
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").Value = "6.613"
$ws.Range("E21").Value = "  +1.36%  "
$ws.Range("D22").Value = "17.40"
$ws.Range("E22").Value = "  -2.25%  "
$ws.Range("D23").Value = "27.207.66"
$ws.Range("E23").Value = "  -1.28%  "
$ws.Range("E24").Value = "  -6.13%  "
$ws.Range("D25").Value = "2.387"
$ws.Range("E25").Value = "  -3.29%  "
$ws.Range("D26").Value = "2.543"
$ws.Range("E26").Value = "  -5.19%  "
$ws.Range("D27").Value = "1.475"
$ws.Range("E27").Value = "  -1.43%  "
$ws.Range("D28").Value = "21.31"
$ws.Range("E28").Value = "  -3.62%  "
$ws.Range("D29").Value = "154.17"
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").Value = "1.987.21"
$ws.Range("E30").Value = "  -1.49%  "
$ws.Range("D31").Value = "134.34"
$ws.Range("E31").Value = "  -1.82%  "
$ws.Range("D32").Value = "4.016"
$ws.Range("E32").Value = "  -1.58%  "
$ws.Range("D33").Value = "6.089"
$ws.Range("E33").Value = "  -4.95%  "
$ws.Range("D34").Value = "0.08738"
$ws.Range("E34").Value = "  -0.73%  "
$ws.Range("E35").Value = "  -4.32%  "
$ws.Range("D36").Value = "1.655"
$ws.Range("E36").Value = "  -3.91%  "
$ws.Range("D37").Value = "0.6917"
$ws.Range("E37").Value = "  -2.61%  "
$ws.Range("D38").Value = "5.443"
$ws.Range("E38").Value = "  -3.71%  "
$ws.Range("D39").Value = "0.2204"
$ws.Range("E39").Value = "  -2.87%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.06350"
$ws.Range("E40").Value = "  -3.12%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "8.805"
$ws.Range("E41").Value = "  -2.51%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.02340"
$ws.Range("E42").Value = "  -3.32%  "
$ws.Range("D43").Value = "1.237"
$ws.Range("E43").Value = "  -4.57%  "
$ws.Range("D44").Value = "14.31"
$ws.Range("E44").Value = "  -4.40%  "
$ws.Range("D45").Value = "0.6498"
$ws.Range("E45").Value = "  -1.95%  "
$ws.Range("D46").Value = "1.002"
$ws.Range("E46").Value = "  +0.43%  "
$ws.Range("D47").Value = "3.849"
$ws.Range("E47").Value = "  -4.81%  "
$ws.Range("D48").Value = "2.144"
$ws.Range("E48").Value = "  -1.98%  "
$ws.Range("D49").Value = "131.05"
$ws.Range("E49").Value = "  -1.51%  "
$ws.Range("D50").Value = "0.07142"
$ws.Range("E50").Value = "  -3.05%  "
$ws.Range("D51").Value = "79.04"
$ws.Range("E51").Value = "  -2.12%  "
